$d = $word.ActiveDocument

# "Sort and Filter products" test case's sub-items are being edited:
#   - "Order products alphabetically and " -> "Order products alphabetically "
#   - "Order products in increasing and decreasing order of price"
#       -> "Order products in order of price"
#   - "Filter by a single product tag at a time" -> "Filter by a single product tag "
# The hidden "_GoBack" bookmark (Word's "last edit location" marker) also moves
# from the end of "Sort and Filter products" to the point in the price bullet
# where text was removed, i.e. right before "order of price".

$d.Content.Find.Execute(
    "Order products alphabetically and ", $true, $false, $false, $false, $false,
    $true, 1, $false, "Order products alphabetically ", 2)

$d.Content.Find.Execute(
    "Order products in increasing and decreasing order of price", $true, $false, $false, $false, $false,
    $true, 1, $false, "Order products in order of price", 2)

$d.Content.Find.Execute(
    "Filter by a single product tag at a time", $true, $false, $false, $false, $false,
    $true, 1, $false, "Filter by a single product tag ", 2)

# Relocate the "_GoBack" bookmark to sit right before "order of price" --
# re-adding a bookmark with the same name moves it rather than duplicating it.
$found = $d.Content
$found.Find.Execute("order of price")
$bmPos = $found.Start
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
